# Update the "dSF" column (F) values for several rows as part of a
# repull/recalculation of the underlying data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    7  = 5
    9  = -2
    10 = -1
    11 = -1
    16 = -3
    19 = 4
    29 = 2
    30 = -1
    36 = -6
    39 = 3
    43 = 0
    45 = -1
    48 = 2
    53 = 1
    60 = 1
    66 = -2
    68 = -1
    72 = 9
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
